$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-51: refreshed crypto price / volume data from the scheduled scraper run.
$ws.Range("D2").Value = "26.106.53"
$ws.Range("E2").Value = "  -4.63%  "
$ws.Range("D3").Value = "1.651.49"
$ws.Range("E3").Value = "  -3.79%  "
$ws.Range("D4").Value = "'1.014"
$ws.Range("E4").Value = "  +0.59%  "
$ws.Range("D5").Value = "'215.83"
$ws.Range("D6").Value = "'0.5106"
$ws.Range("E6").Value = "  -3.33%  "
$ws.Range("D7").Value = "'1.015"
$ws.Range("E7").Value = "  +0.69%  "
$ws.Range("E8").Value = "  -2.34%  "
$ws.Range("E9").Value = "  -3.59%  "
$ws.Range("D10").Value = "'19.64"
$ws.Range("E10").Value = "  -5.81%  "
$ws.Range("D11").Value = "'0.07774"
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("D12").Value = "1.655.66"
$ws.Range("E12").Value = "  -3.78%  "
$ws.Range("D13").Value = "'4.272"
$ws.Range("E13").Value = "  -4.59%  "
$ws.Range("D14").Value = "1.881.08"
$ws.Range("E14").Value = "  -3.65%  "
$ws.Range("D15").Value = "'0.5476"
$ws.Range("E15").Value = "  -5.67%  "
$ws.Range("D16").Value = "0.0₅7980"
$ws.Range("E16").Value = "  -2.55%  "
$ws.Range("D17").Value = "'63.74"
$ws.Range("E17").Value = "  -6.05%  "
$ws.Range("D18").Value = "26.122.65"
$ws.Range("E18").Value = "  -4.60%  "
$ws.Range("D19").Value = "'1.013"
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("D20").Value = "'207.19"
$ws.Range("E20").Value = "  -5.65%  "
$ws.Range("D21").Value = "'4.394"
$ws.Range("E21").Value = "  -5.69%  "
$ws.Range("D22").Value = "'10.07"
$ws.Range("E22").Value = "  -3.47%  "
$ws.Range("D23").Value = "'6.032"
$ws.Range("E23").Value = "  -0.40%  "
$ws.Range("D24").Value = "'1.015"
$ws.Range("E24").Value = "  +0.68%  "
$ws.Range("E25").Value = "  +7.50%  "
$ws.Range("D26").Value = "'144.16"
$ws.Range("E26").Value = "  -0.77%  "
$ws.Range("D27").Value = "'0.1169"
$ws.Range("E27").Value = "  -3.21%  "
$ws.Range("E28").Value = "  -4.03%  "
$ws.Range("D29").Value = "'15.79"
$ws.Range("E29").Value = "  -2.60%  "
$ws.Range("D30").Value = "'0.05096"
$ws.Range("E30").Value = "  -4.81%  "
$ws.Range("D31").Value = "'1.244"
$ws.Range("E31").Value = "  -3.93%  "
$ws.Range("D32").Value = "'3.341"
$ws.Range("E32").Value = "  -4.16%  "
$ws.Range("E33").Value = "  -4.85%  "
$ws.Range("D34").Value = "'1.545"
$ws.Range("E34").Value = "  -5.94%  "
$ws.Range("D35").Value = "'2.367"
$ws.Range("E35").Value = "  -1.38%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'0.9165"
$ws.Range("E36").Value = "  -4.04%  "
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").Value = "'2.679"
$ws.Range("E37").Value = "  -5.64%  "
$ws.Range("D38").Value = "1.173.08"
$ws.Range("E38").Value = "  -0.91%  "
$ws.Range("D39").Value = "'0.5684"
$ws.Range("E39").Value = "  -3.44%  "
$ws.Range("D40").Value = "'0.01578"
$ws.Range("E40").Value = "  -4.56%  "
$ws.Range("D41").Value = "'2.586"
$ws.Range("E41").Value = "  +0.44%  "
$ws.Range("D42").Value = "'1.016"
$ws.Range("E42").Value = "  +0.76%  "
$ws.Range("D43").Value = "'5.675"
$ws.Range("E43").Value = "  -2.54%  "
$ws.Range("D44").Value = "'0.8247"
$ws.Range("E44").Value = "  -1.94%  "
$ws.Range("D45").Value = "'100.31"
$ws.Range("E45").Value = "  -0.91%  "
$ws.Range("D46").Value = "1.792.64"
$ws.Range("E46").Value = "  -3.58%  "
$ws.Range("D47").Value = "0.0₈112"
$ws.Range("E47").Value = "  -5.21%  "
$ws.Range("D48").Value = "'0.4555"
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("D49").Value = "'1.010"
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("D50").Value = "'55.24"
$ws.Range("E50").Value = "  -4.09%  "
$ws.Range("D51").Value = "'7.836"
$ws.Range("E51").Value = "  -4.16%  "
